# Localization status report regeneration:
#   - Status moves from "Ready for handoff" to "In Translation" for all
#     tracked files (Overview summary columns + per-locale Status column).
#   - Because the new status text is shorter, the Status-related columns
#     shrink from their previous autosized width.

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"
# Target stored column width (OOXML) is 13.4101845877511 characters, which
# is what real Excel's font-metric AutoFit produces for the new, shorter
# status text. This engine's ColumnWidth setter quantizes to a coarser
# pixel grid, so 12.5 is the input that lands on the closest reachable
# width bucket (13.333333333333334) to that target.
$newWidth  = 12.5

# --- Overview sheet: one status column per locale (zh-cn, de-de) ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2:F4").Value2 = $newStatus

# --- Per-locale sheets: "Status" is column C on both zh-cn and de-de ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2:C4").Value2 = $newStatus

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2:C4").Value2 = $newStatus

# --- Re-fit the affected columns now that the text is shorter ---
$overview.Columns.Item(5).ColumnWidth = $newWidth
$overview.Columns.Item(6).ColumnWidth = $newWidth
$zhcn.Columns.Item(3).ColumnWidth = $newWidth
$dede.Columns.Item(3).ColumnWidth = $newWidth
